$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.893.17'
$ws.Range('E2').Value = '  +2.09%  '
$ws.Range('D3').Value = '1.811.23'
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Value = '''313.73'
$ws.Range('E5').Value = '  +3.24%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').Value = '''0.4289'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').Value = '''0.3697'
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('D9').Value = '''0.07258'
$ws.Range('E9').Value = '  +3.04%  '
$ws.Range('D10').Value = '''0.8654'
$ws.Range('E10').Value = '  +4.25%  '
$ws.Range('D11').Value = '2.039.40'
$ws.Range('E11').Value = '  +16.50%  '
$ws.Range('D12').Value = '''21.29'
$ws.Range('E12').Value = '  +5.72%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '''6.627'
$ws.Range('E13').Value = '  +3.98%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''5.399'
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('D15').Value = '''0.06933'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('D16').Value = '''80.76'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = '''1.012'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').Value = '''0.000008904'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('D19').Value = '''1.006'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('D21').Value = '26.936.85'
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').Value = '''5.201'
$ws.Range('E22').Value = '  +4.20%  '
$ws.Range('D23').Value = '''10.95'
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('D24').Value = '2.278.04'
$ws.Range('E24').Value = '  +15.66%  '
$ws.Range('D25').Value = '''154.24'
$ws.Range('E25').Value = '  +1.66%  '
$ws.Range('D26').Value = '''1.884'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').Value = '''5.245'
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('D29').Value = '''1.924'
$ws.Range('E29').Value = '  +15.04%  '
$ws.Range('D30').Value = '''114.72'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('D32').Value = '''0.7440'
$ws.Range('E32').Value = '  +2.84%  '
$ws.Range('E33').Value = '  +4.42%  '
$ws.Range('D34').Value = '''4.433'
$ws.Range('E34').Value = '  +3.00%  '
$ws.Range('D35').Value = '''2.805'
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('D36').Value = '''1.006'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  +5.50%  '
$ws.Range('D38').Value = '''0.05232'
$ws.Range('E38').Value = '  +2.77%  '
$ws.Range('D39').Value = '''0.01928'
$ws.Range('E39').Value = '  +2.69%  '
$ws.Range('D40').Value = '''0.5092'
$ws.Range('E40').Value = '  +4.15%  '
$ws.Range('D41').Value = '''2.744'
$ws.Range('E41').Value = '  +9.99%  '
$ws.Range('D42').Value = '''0.1650'
$ws.Range('E42').Value = '  +2.98%  '
$ws.Range('D43').Value = '''6.510'
$ws.Range('E43').Value = '  +5.77%  '
$ws.Range('D44').Value = '''8.317'
$ws.Range('E44').Value = '  +4.17%  '
$ws.Range('D45').Value = '''107.38'
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('D46').Value = '''10.40'
$ws.Range('E46').Value = '  +3.71%  '
$ws.Range('D47').Value = '''1.006'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.650'
$ws.Range('E48').Value = '  +5.46%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = '''0.4556'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').Value = '''0.06267'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').Value = '''1.800'
$ws.Range('E51').Value = '  +5.16%  '
